# Set BB-specific price trends.
# Updates the "N" column inputs (near-term annual growth trend) on the
# "Prices" sheet for each of the five price-level rows (Investment,
# Consumption, Imports, Exports, GDP), which in turn changes the
# dependent formulas (N4, E5:M9) via recalculation. Also bumps the
# number format on those input cells from 2 to 3 decimal places and
# moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

# New BB-specific price trend values (column N, rows 5-9)
$ws.Range("N5").Value = 0.06
$ws.Range("N6").Value = 0.006
$ws.Range("N7").Value = 0.01
$ws.Range("N8").Value = 0.015
$ws.Range("N9").Value = 0.025

# These cells (and the dependent N4 total) now display three decimals
# instead of two.
$ws.Range("N4:N9").NumberFormat = "0.000"

# Force a full recalculation so the formulas in N4 and E5:M9 pick up
# the new trend values.
$excel.CalculateFull() | Out-Null

# Update the active selection on the sheet to match the saved view.
$ws.Activate() | Out-Null
$ws.Range("N10").Select() | Out-Null
